$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 213-214, pushing the existing data (old rows
# 213-256) down to 215-258.
$ws.Rows("213:214").Insert()

# New row 213
$ws.Range("A213").Value = 7
$ws.Range("B213").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C213").Value = 'Ñuble'
$ws.Range("D213").Value = 44785
$ws.Range("D213").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E213").Value = 16
$ws.Range("F213").Value = 100112006
$ws.Range("G213").Value = 'Repollo'
$ws.Range("H213").Value = 'Crespo record'
$ws.Range("I213").Value = 'Primera'
$ws.Range("J213").Value = 200
$ws.Range("K213").Value = 1100
$ws.Range("L213").Value = 1300
$ws.Range("M213").Value = 1200
$ws.Range("N213").Value = '$/unidad'
$ws.Range("O213").Value = 'Provincia de Diguillín'
$ws.Range("P213").Value = 1200
$ws.Range("Q213").Value = 1
$ws.Range("R213").Value = 'Hortaliza'

# New row 214
$ws.Range("A214").Value = 7
$ws.Range("B214").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C214").Value = 'Ñuble'
$ws.Range("D214").Value = 44785
$ws.Range("D214").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E214").Value = 16
$ws.Range("F214").Value = 100112006
$ws.Range("G214").Value = 'Repollo'
$ws.Range("H214").Value = 'Crespo record'
$ws.Range("I214").Value = 'Segunda'
$ws.Range("J214").Value = 200
$ws.Range("K214").Value = 900
$ws.Range("L214").Value = 900
$ws.Range("M214").Value = 900
$ws.Range("N214").Value = '$/unidad'
$ws.Range("O214").Value = 'Provincia de Diguillín'
$ws.Range("P214").Value = 900
$ws.Range("Q214").Value = 1
$ws.Range("R214").Value = 'Hortaliza'
